$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D (Price) and E (Volume(1h)) columns stay text so values like "1.000" or
# "0.9994" are not silently coerced into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.498.08"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.912.18"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "244.42"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "0.4850"
$ws.Range("E7").Value = "  +3.31%  "
$ws.Range("D8").Value = "0.2897"
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("D9").Value = "0.06711"
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("D10").Value = "109.61"
$ws.Range("E10").Value = "  +3.04%  "
$ws.Range("D11").Value = "19.37"
$ws.Range("E11").Value = "  +7.31%  "
$ws.Range("D12").Value = "1.912.42"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "0.07531"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").Value = "5.261"
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").Value = "0.6666"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "271.82"
$ws.Range("E16").Value = "  -4.95%  "
$ws.Range("D17").Value = "30.484.46"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "0.9994"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "0.000007533"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").Value = "12.86"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").Value = "2.165.53"
$ws.Range("D22").Value = "5.498"
$ws.Range("E22").Value = "  +5.44%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "6.399"
$ws.Range("E24").Value = "  +3.47%  "
$ws.Range("D25").Value = "9.410"
$ws.Range("E25").Value = "  +1.94%  "
$ws.Range("D26").Value = "163.66"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("D27").Value = "20.18"
$ws.Range("E27").Value = "  -4.49%  "
$ws.Range("D28").Value = "2.098"
$ws.Range("E28").Value = "  +3.19%  "
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").Value = "1.399"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").Value = "4.120"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").Value = "4.037"
$ws.Range("D33").Value = "0.04981"
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").Value = "0.7275"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "1.134"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").Value = "0.9999"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "0.02029"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("D39").Value = "2.666"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "111.02"
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("D41").Value = "2.013"
$ws.Range("E41").Value = "  -1.71%  "
$ws.Range("D42").Value = "0.4415"
$ws.Range("E42").Value = "  +5.26%  "
$ws.Range("D43").Value = "0.8668"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").Value = "5.844"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").Value = "0.9989"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "67.76"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").Value = "7.406"
$ws.Range("E47").Value = "  +4.07%  "
$ws.Range("D48").Value = "9.229"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").Value = "0.1243"
$ws.Range("E49").Value = "  +3.41%  "
$ws.Range("D50").Value = "47.37"
$ws.Range("E50").Value = "  -9.57%  "
$ws.Range("D51").Value = "1.465"
$ws.Range("E51").Value = "  +6.64%  "
